$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sub-Tasks")

$subtaskData = @(
  @("Subtask", "Description"),
  @("Conference planning", "Involves logistical coordination, meetings, ordering items, maintaining invite lists, answering questions, accepting RSVPs, Putting together agendas, writing , creating and Implementing promotional strategy associated with conference, etc."),
  @("At event", "Involves being on site at event location supervising logistics and posting to social media if appropriate for event. "),
  @("Invitations", "Involves creating, editing, emailing, printing and delivering of various event/conference invitations. "),
  @("Update", "Involves updating website, employee directories, "),
  @("Creation", ""),
  @("RFP submission", "Involves assisting center directors/group leaders in RFP delivery to win event or conference bid. "),
  @("Graphics design", ""),
  @("Set up & instruction", ""),
  @("Small graphic edit", "Involves editing, revising, incorporating feedback into a small scale/size graphic"),
  @("Medium graphic edit", "Involves editing, revising, incorporating feedback into a moderate scale/size graphic"),
  @("Large graphic edit", "Involves editing, revising, incorporating feedback into a large scale/size graphic"),
  @("Small graphic creation", "Involves researching/gathering information, ideation and production of a small scale/size graphic"),
  @("Medium graphic creation", "Involves researching/gathering information, ideation and production of a moderate scale/size graphic"),
  @("Large graphic creation", "Involves researching/gathering information, ideation and production of a large scale/size graphic"),
  @("Animation", "Involves actually setting graphic components into motion and cultivating weight, feel and timing "),
  @("Photography - Preproduction", "Involves planning the photoshoot: gathering necessary materials and personnel."),
  @("Photography - Production", "Involves executing the photoshoot"),
  @("Photography - Postproduction", "Involves file upload, color correction, exposure balancing, general digital touch-ups"),
  @("Dealing with ensemble", ""),
  @("File/Folder cleanup", "Involves cleaning up files and reorganizing the workspace"),
  @("Printing", "Involves printing of various graphic materials "),
  @("Communication", ""),
  @("Training", ""),
  @("File transfer", "Involves transferring of data from camera to VTTI server or from VTTI server to a media contact/YouTube etc.  "),
  @("Answering questions", ""),
  @("Demo Site", ""),
  @("Hardware", ""),
  @("Creating content", "Involves creation of content by writing and may include some graphic manipulation,"),
  @("Writing content", "Involves writing for social media, news releases, for website, for blog posts, for conference and event materials"),
  @("Editing content", "Involves modifying or editing existing content such as news releases, web content, blog posts, and  graphic materials"),
  @("Ordering", "Involves filling out PAF and working with fiscal team to get product/service from purchase to received stages "),
  @("Coordination", "Involves attendance  and coordination for Tours/Meetings/Events/Interviews/Demos. May involve attending meetings, examining researchers calendars, sending emails, submitting necessary forms and answering phone calls to setup activity."),
  @("Updating Invitation List", "Involves updating, adding and subtracting contact names from various event invitation lists as well as the VTTI Master list. "),
  @("Training", "Involves signing up and attending various trainings either mandatory or for professional development. "),
  @("Hiring ", "Involves vetting candidates for employment or internships. Can involve resume review, interviews, and any paperwork/training involved to support new team member. "),
  @("Video - Preproduction", "Involves planning principal photography, scheduling, gathering materials and storyboarding"),
  @("Video - Production", "Involves executing principal photography and audio"),
  @("Video - Postproduction", "Involves using final cut pro, after effects  or adobe premier to edit and deliver video requests. "),
)

# Clear the existing Subtask/Description table before rewriting it with the
# updated set of rows (several subtasks gained descriptions, the old
# "Filming"/"Video editing"/"Photo editing" rows were replaced by a proper
# Photography/Video pre-post-production breakdown, and a "Small graphic
# creation" + extra "Training" row were added).
$ws.Range("A1:B35").ClearContents()

$r = 1
foreach ($row in $subtaskData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne "") {
        $ws.Cells.Item($r, 2).Value = $row[1]
    }
    $r = $r + 1
}

# Highlight the still-undocumented "Communication" and "Answering questions"
# subtasks in yellow so they stand out as needing a description later.
$ws.Range("A23").Interior.Color = 65535
$ws.Range("A26").Interior.Color = 65535
